$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2,2).Value = 0.2403245088184178
$ws.Cells.Item(2,3).Value = 0.0634223395761353
$ws.Cells.Item(2,5).Value = 0.177508499017442
$ws.Cells.Item(2,6).Value = 0.4443680307746121
$ws.Cells.Item(2,7).Value = 0.3486100814986131
$ws.Cells.Item(2,8).Value = 0.5268670116760177
$ws.Cells.Item(2,11).Value = 0.2275448627089958
$ws.Cells.Item(2,13).Value = 0.2035923660319341
$ws.Cells.Item(2,14).Value = 1.258827391824198
$ws.Cells.Item(2,15).Value = 1.687441793970777
$ws.Cells.Item(3,2).Value = 0.2108832419778537
$ws.Cells.Item(3,3).Value = 0.06128855930615629
$ws.Cells.Item(3,5).Value = 0.1659538764275936
$ws.Cells.Item(3,6).Value = 0.3878228170618172
$ws.Cells.Item(3,7).Value = 0.351215470347384
$ws.Cells.Item(3,8).Value = 0.5312230593745824
$ws.Cells.Item(3,11).Value = 0.1987109884060345
$ws.Cells.Item(3,13).Value = 0.1821749732985012
$ws.Cells.Item(3,14).Value = 1.271542974129794
$ws.Cells.Item(3,15).Value = 1.701875517631947
$ws.Cells.Item(4,2).Value = 0.1927774303031526
$ws.Cells.Item(4,3).Value = 0.05997247440143383
$ws.Cells.Item(4,5).Value = 0.1589817332416956
$ws.Cells.Item(4,6).Value = 0.3531389305169483
$ws.Cells.Item(4,7).Value = 0.3530736056629351
$ws.Cells.Item(4,8).Value = 0.5341203111805584
$ws.Cells.Item(4,11).Value = 0.1809326983169655
$ws.Cells.Item(4,13).Value = 0.1690761615674816
$ws.Cells.Item(4,14).Value = 1.27979691885438
$ws.Cells.Item(4,15).Value = 1.711744608653845
$ws.Cells.Item(5,2).Value = 0.1853923994532352
$ws.Cells.Item(5,3).Value = 0.05943470958581543
$ws.Cells.Item(5,5).Value = 0.1561711622769906
$ws.Cells.Item(5,6).Value = 0.3390132514313251
$ws.Cells.Item(5,7).Value = 0.3538957139314221
$ws.Cells.Item(5,8).Value = 0.5353569737162402
$ws.Cells.Item(5,11).Value = 0.1736696910071203
$ws.Cells.Item(5,13).Value = 0.16375131039878
$ws.Cells.Item(5,14).Value = 1.283272786664877
$ws.Cells.Item(5,15).Value = 1.716019384014231
$ws.Cells.Item(6,2).Value = 0.184165725771237
$ws.Cells.Item(6,3).Value = 0.05934532786395152
$ws.Cells.Item(6,5).Value = 0.1557063147960278
$ws.Cells.Item(6,6).Value = 0.336668177824194
$ws.Cells.Item(6,7).Value = 0.3540361420864357
$ws.Cells.Item(6,8).Value = 0.5355657042473823
$ws.Cells.Item(6,11).Value = 0.1724625873693668
$ws.Cells.Item(6,13).Value = 0.1628679110558267
$ws.Cells.Item(6,14).Value = 1.283856736707737
$ws.Cells.Item(6,15).Value = 1.71674448705997
$ws.Cells.Item(7,2).Value = 0.1926778599720649
$ws.Cells.Item(7,3).Value = 0.05996522772888824
$ws.Cells.Item(7,5).Value = 0.1589437050678342
$ws.Cells.Item(7,6).Value = 0.3529483938344953
$ws.Cells.Item(7,7).Value = 0.3530844302157661
$ws.Cells.Item(7,8).Value = 0.5341367624275293
$ws.Cells.Item(7,11).Value = 0.1808348200048613
$ws.Cells.Item(7,13).Value = 0.1690042959648963
$ws.Cells.Item(7,14).Value = 1.279843340802312
$ws.Cells.Item(7,15).Value = 1.711801235362017
$ws.Cells.Item(8,2).Value = 0.2301793981734477
$ws.Cells.Item(8,3).Value = 0.06268786907369872
$ws.Cells.Item(8,5).Value = 0.173498955447144
$ws.Cells.Item(8,6).Value = 0.4248636149813478
$ws.Cells.Item(8,7).Value = 0.349454738216366
$ws.Cells.Item(8,8).Value = 0.5283227894467046
$ws.Cells.Item(8,11).Value = 0.2176185944867086
$ws.Cells.Item(8,13).Value = 0.1961969768694019
$ws.Cells.Item(8,14).Value = 1.26311907040542
$ws.Cells.Item(8,15).Value = 1.692209529658101
$ws.Cells.Item(9,2).Value = 0.3034746813862625
$ws.Cells.Item(9,3).Value = 0.06797822700022493
$ws.Cells.Item(9,5).Value = 0.2030218565510253
$ws.Cells.Item(9,6).Value = 0.5661985755041457
$ws.Cells.Item(9,7).Value = 0.3443910777095738
$ws.Cells.Item(9,8).Value = 0.5186867381479487
$ws.Cells.Item(9,11).Value = 0.2891467992538992
$ws.Cells.Item(9,13).Value = 0.2499310250605049
$ws.Cells.Item(9,14).Value = 1.233863866737092
$ws.Cells.Item(9,15).Value = 1.661783676680798
$ws.Cells.Item(10,2).Value = 0.3571576059917447
$ws.Cells.Item(10,3).Value = 0.0718334999612722
$ws.Cells.Item(10,5).Value = 0.2253239855342386
$ws.Cells.Item(10,6).Value = 0.6702781546542269
$ws.Cells.Item(10,7).Value = 0.3419285818229625
$ws.Cells.Item(10,8).Value = 0.5126815210477389
$ws.Cells.Item(10,11).Value = 0.3413131436348351
$ws.Cells.Item(10,13).Value = 0.2896633022169723
$ws.Cells.Item(10,14).Value = 1.214525623300304
$ws.Cells.Item(10,15).Value = 1.644311272702822
$ws.Cells.Item(11,2).Value = 0.3815396172784631
$ws.Cells.Item(11,3).Value = 0.0735801437845538
$ws.Cells.Item(11,5).Value = 0.2356057813172754
$ws.Cells.Item(11,6).Value = 0.7176906081379002
$ws.Cells.Item(11,7).Value = 0.341082646327358
$ws.Cells.Item(11,8).Value = 0.5101825421181729
$ws.Cells.Item(11,11).Value = 0.3649579948809958
$ws.Cells.Item(11,13).Value = 0.3077948295141013
$ws.Cells.Item(11,14).Value = 1.206195565599351
$ws.Cells.Item(11,15).Value = 1.637424486698478
$ws.Cells.Item(12,2).Value = 0.3907664881040205
$ws.Cells.Item(12,3).Value = 0.07424048683675721
$ws.Cells.Item(12,5).Value = 0.2395190524741935
$ws.Cells.Item(12,6).Value = 0.7356546913071611
$ws.Cells.Item(12,7).Value = 0.3408018440481655
$ws.Cells.Item(12,8).Value = 0.5092696994810879
$ws.Cells.Item(12,11).Value = 0.3738989468864986
$ws.Cells.Item(12,13).Value = 0.3146689789071573
$ws.Cells.Item(12,14).Value = 1.203108318842702
$ws.Cells.Item(12,15).Value = 1.634969432203391
$ws.Cells.Item(13,2).Value = 0.3887795944305026
$ws.Cells.Item(13,3).Value = 0.0740983186080939
$ws.Cells.Item(13,5).Value = 0.2386753773722177
$ws.Cells.Item(13,6).Value = 0.7317853510981394
$ws.Cells.Item(13,7).Value = 0.3408605599856287
$ws.Cells.Item(13,8).Value = 0.5094648082578388
$ws.Cells.Item(13,11).Value = 0.371973931767684
$ws.Cells.Item(13,13).Value = 0.3131881473811902
$ws.Cells.Item(13,14).Value = 1.203770225284735
$ws.Cells.Item(13,15).Value = 1.635491372466888
$ws.Cells.Item(14,2).Value = 0.3822988420444347
$ws.Cells.Item(14,3).Value = 0.07363449237503517
$ws.Cells.Item(14,5).Value = 0.235927331206085
$ws.Cells.Item(14,6).Value = 0.7191683204515869
$ws.Cells.Item(14,7).Value = 0.3410587516641499
$ws.Cells.Item(14,8).Value = 0.5101067713638443
$ws.Cells.Item(14,11).Value = 0.3656938325248404
$ws.Cells.Item(14,13).Value = 0.3083602073742924
$ws.Cells.Item(14,14).Value = 1.205940229790308
$ws.Cells.Item(14,15).Value = 1.637219444082262
$ws.Cells.Item(15,2).Value = 0.3783283886293418
$ws.Cells.Item(15,3).Value = 0.07335024429593773
$ws.Cells.Item(15,5).Value = 0.2342466539514163
$ws.Cells.Item(15,6).Value = 0.7114413442032514
$ws.Cells.Item(15,7).Value = 0.3411853012699524
$ws.Cells.Item(15,8).Value = 0.5105043502280324
$ws.Cells.Item(15,11).Value = 0.3618454037045069
$ws.Cells.Item(15,13).Value = 0.3054040115566181
$ws.Cells.Item(15,14).Value = 1.207278167461194
$ws.Cells.Item(15,15).Value = 1.638297845088118
$ws.Cells.Item(16,2).Value = 0.3555633631234514
$ws.Cells.Item(16,3).Value = 0.07171920504647744
$ws.Cells.Item(16,5).Value = 0.2246548051200463
$ws.Cells.Item(16,6).Value = 0.6671810134426437
$ws.Cells.Item(16,7).Value = 0.3419893928534421
$ws.Cells.Item(16,8).Value = 0.512849519068304
$ws.Cells.Item(16,11).Value = 0.3397661269289642
$ws.Cells.Item(16,13).Value = 0.2884795067953405
$ws.Cells.Item(16,14).Value = 1.215079410354953
$ws.Cells.Item(16,15).Value = 1.644782715805704
$ws.Cells.Item(17,2).Value = 0.3415874999280959
$ws.Cells.Item(17,3).Value = 0.07071675368651142
$ws.Cells.Item(17,5).Value = 0.2188055758941516
$ws.Cells.Item(17,6).Value = 0.6400460337125793
$ws.Cells.Item(17,7).Value = 0.3425529916392662
$ws.Cells.Item(17,8).Value = 0.5143478264894199
$ws.Cells.Item(17,11).Value = 0.3261988606507487
$ws.Cells.Item(17,13).Value = 0.2781114436401424
$ws.Cells.Item(17,14).Value = 1.219984838334618
$ws.Cells.Item(17,15).Value = 1.649032977244673
$ws.Cells.Item(18,2).Value = 0.3335453486667461
$ws.Cells.Item(18,3).Value = 0.07013950097208976
$ws.Cells.Item(18,5).Value = 0.2154540874860018
$ws.Cells.Item(18,6).Value = 0.6244449056556647
$ws.Cells.Item(18,7).Value = 0.3429029682884561
$ws.Cells.Item(18,8).Value = 0.5152315286130502
$ws.Cells.Item(18,11).Value = 0.3183872737548654
$ws.Cells.Item(18,13).Value = 0.2721533887230336
$ws.Cells.Item(18,14).Value = 1.222850266828594
$ws.Cells.Item(18,15).Value = 1.651577517741615
$ws.Cells.Item(19,2).Value = 0.3308218110394137
$ws.Cells.Item(19,3).Value = 0.0699439396013446
$ws.Cells.Item(19,5).Value = 0.2143215304530344
$ws.Cells.Item(19,6).Value = 0.619163680173358
$ws.Cells.Item(19,7).Value = 0.3430258942472264
$ws.Cells.Item(19,8).Value = 0.5155344993604274
$ws.Cells.Item(19,11).Value = 0.3157410362914845
$ws.Cells.Item(19,13).Value = 0.2701370202601652
$ws.Cells.Item(19,14).Value = 1.223828000901825
$ws.Cells.Item(19,15).Value = 1.652456208302027
$ws.Cells.Item(20,2).Value = 0.343075631190203
$ws.Cells.Item(20,3).Value = 0.07082353594215363
$ws.Cells.Item(20,5).Value = 0.2194269069646282
$ws.Cells.Item(20,6).Value = 0.642933953830422
$ws.Cells.Item(20,7).Value = 0.342490323699117
$ws.Cells.Item(20,8).Value = 0.5141860610502391
$ws.Cells.Item(20,11).Value = 0.3276439564376403
$ws.Cells.Item(20,13).Value = 0.2792145853451586
$ws.Cells.Item(20,14).Value = 1.219458097795339
$ws.Cells.Item(20,15).Value = 1.648570188815896
$ws.Cells.Item(21,2).Value = 0.3842025648892218
$ws.Cells.Item(21,3).Value = 0.07377075887582407
$ws.Cells.Item(21,5).Value = 0.2367339609521224
$ws.Cells.Item(21,6).Value = 0.7228739723491628
$ws.Cells.Item(21,7).Value = 0.3409994642298884
$ws.Cells.Item(21,8).Value = 0.5099173030743174
$ws.Cells.Item(21,11).Value = 0.3675388022860773
$ws.Cells.Item(21,13).Value = 0.3097780702846009
$ws.Cells.Item(21,14).Value = 1.205301023923951
$ws.Cells.Item(21,15).Value = 1.63670771858861
$ws.Cells.Item(22,2).Value = 0.4110457877701776
$ws.Cells.Item(22,3).Value = 0.07569066339230801
$ws.Cells.Item(22,5).Value = 0.2481604763730445
$ws.Cells.Item(22,6).Value = 0.7751780083420101
$ws.Cells.Item(22,7).Value = 0.3402555865168893
$ws.Cells.Item(22,8).Value = 0.5073224849838667
$ws.Cells.Item(22,11).Value = 0.3935373214889069
$ws.Cells.Item(22,13).Value = 0.3298004450258532
$ws.Cells.Item(22,14).Value = 1.196440055029207
$ws.Cells.Item(22,15).Value = 1.629845754397763
$ws.Cells.Item(23,2).Value = 0.3967224823425397
$ws.Cells.Item(23,3).Value = 0.07466656378663572
$ws.Cells.Item(23,5).Value = 0.2420513232286936
$ws.Cells.Item(23,6).Value = 0.7472568307830727
$ws.Cells.Item(23,7).Value = 0.3406314870881317
$ws.Cells.Item(23,8).Value = 0.5086895444183597
$ws.Cells.Item(23,11).Value = 0.37966845330007
$ws.Cells.Item(23,13).Value = 0.3191098190133346
$ws.Cells.Item(23,14).Value = 1.201133494836885
$ws.Cells.Item(23,15).Value = 1.633426542188829
$ws.Cells.Item(24,2).Value = 0.3424028695636991
$ws.Cells.Item(24,3).Value = 0.07077526257903344
$ws.Cells.Item(24,5).Value = 0.2191459679676413
$ws.Cells.Item(24,6).Value = 0.6416283278902171
$ws.Cells.Item(24,7).Value = 0.3425185750303754
$ws.Cells.Item(24,8).Value = 0.5142591257243083
$ws.Cells.Item(24,11).Value = 0.3269906646805794
$ws.Cells.Item(24,13).Value = 0.2787158466201163
$ws.Cells.Item(24,14).Value = 1.219696096296058
$ws.Cells.Item(24,15).Value = 1.648779100842901
$ws.Cells.Item(25,2).Value = 0.2836744358168346
$ws.Cells.Item(25,3).Value = 0.06655245339403848
$ws.Cells.Item(25,5).Value = 0.1949287365514607
$ws.Cells.Item(25,6).Value = 0.5279251897347166
$ws.Cells.Item(25,7).Value = 0.3455404497146546
$ws.Cells.Item(25,8).Value = 0.5211047508732989
$ws.Cells.Item(25,11).Value = 0.2698629969668787
$ws.Cells.Item(25,13).Value = 0.2353502696692118
$ws.Cells.Item(25,14).Value = 1.241399417789648
$ws.Cells.Item(25,15).Value = 1.669158049056648
